$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "개설학과" (department) column values for rows 2-9, keeping the cell style.
$ws.Range("B2:B9").ClearContents()

# Clear the "교과구분" (subject category) value for row 2 only.
$ws.Range("D2").ClearContents()

# Update selection to B9 to match the saved view state.
$ws.Range("B9").Select()
